# Apply cryptocurrency price/volume updates from the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Prefix with an apostrophe so Excel always stores the value as literal
    # text (never auto-converted to a number/date), then strip the resulting
    # quote-prefix style so cell formatting matches the original (unstyled) cells.
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

Set-TextValue "D2" "46.916.87"
Set-TextValue "E2" "  +5.79%  "
Set-TextValue "D3" "2.308.65"
Set-TextValue "E3" "  +3.70%  "
Set-TextValue "E4" "  -0.71%  "
Set-TextValue "D5" "304.09"
Set-TextValue "E5" "  +1.43%  "
Set-TextValue "D6" "102.03"
Set-TextValue "E6" "  +13.11%  "
Set-TextValue "D7" "0.570"
Set-TextValue "E7" "  +1.47%  "
Set-TextValue "E9" "  +7.86%  "
Set-TextValue "D10" "37.04"
Set-TextValue "E10" "  +12.06%  "
Set-TextValue "D11" "0.0804"
Set-TextValue "E11" "  +2.56%  "
Set-TextValue "D12" "7.47"
Set-TextValue "E12" "  +7.31%  "
Set-TextValue "E13" "  +0.23%  "
Set-TextValue "D14" "2.658.50"
Set-TextValue "E14" "  +3.51%  "
Set-TextValue "D15" "2.302.84"
Set-TextValue "E15" "  +3.68%  "
Set-TextValue "D16" "14.06"
Set-TextValue "E16" "  +4.50%  "
Set-TextValue "D17" "0.824"
Set-TextValue "E17" "  +5.64%  "
Set-TextValue "D18" "46.894.26"
Set-TextValue "E18" "  +6.21%  "
Set-TextValue "D19" "13.75"
Set-TextValue "E19" "  +24.55%  "
Set-TextValue "D20" "0.0₃0951"
Set-TextValue "E20" "  +4.19%  "
Set-TextValue "D21" "6.14"
Set-TextValue "E21" "  +3.61%  "
Set-TextValue "D22" "67.00"
Set-TextValue "E22" "  +3.54%  "
Set-TextValue "D23" "249.23"
Set-TextValue "E23" "  +4.26%  "
Set-TextValue "D24" "2.96"
Set-TextValue "E24" "  +5.29%  "
Set-TextValue "D25" "1.98"
Set-TextValue "E25" "  +6.17%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -1.25%  "
Set-TextValue "D27" "44.37"
Set-TextValue "E27" "  +15.67%  "
Set-TextValue "D28" "2.29"
Set-TextValue "E28" "  +2.23%  "
Set-TextValue "D29" "9.97"
Set-TextValue "E29" "  +6.39%  "
Set-TextValue "D30" "20.27"
Set-TextValue "E30" "  +4.05%  "
Set-TextValue "D31" "5.83"
Set-TextValue "E31" "  +7.89%  "
Set-TextValue "B32" "Hedera"
Set-TextValue "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.0811"
Set-TextValue "E32" "  +8.09%  "
Set-TextValue "B33" "Monero"
Set-TextValue "C33" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D33" "147.78"
Set-TextValue "E33" "  -0.53%  "
Set-TextValue "D34" "2.64"
Set-TextValue "E34" "  +3.59%  "
Set-TextValue "E35" "  +12.78%  "
Set-TextValue "D36" "0.115"
Set-TextValue "E36" "  +11.37%  "
Set-TextValue "E37" "  +3.14%  "
Set-TextValue "E38" "  +7.53%  "
Set-TextValue "D39" "16.30"
Set-TextValue "E39" "  +23.40%  "
Set-TextValue "D40" "4.05"
Set-TextValue "E40" "  +14.56%  "
Set-TextValue "D41" "3.48"
Set-TextValue "E41" "  +8.97%  "
Set-TextValue "E42" "  +1.19%  "
Set-TextValue "D43" "2.02"
Set-TextValue "E43" "  +13.48%  "
Set-TextValue "E44" "  -0.69%  "
Set-TextValue "D45" "1.866.02"
Set-TextValue "E45" "  +2.85%  "
Set-TextValue "D46" "88.86"
Set-TextValue "E46" "  +19.73%  "
Set-TextValue "E47" "  +10.39%  "
Set-TextValue "D48" "75.13"
Set-TextValue "E48" "  +12.42%  "
Set-TextValue "D49" "4.92"
Set-TextValue "E49" "  +10.80%  "
Set-TextValue "D50" "97.49"
Set-TextValue "E50" "  +3.45%  "
Set-TextValue "E51" "  +6.02%  "
